$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "VALOR MORA" total and "Cant. Trabajadores" count
$ws.Range("E11").Value = 568534
$ws.Range("C13").Value = 1

# 2. Re-sort worker 1 (SHIRLY NANY MURFORD PUERTA / CC 23243902) period rows (16-21)
#    into ascending chronological order, keeping each period's own "Valor Mora" value attached.
$periods = @(2405, 2407, 2408, 2409, 2410, 2411)
$values  = @(104000, 104000, 104000, 104000, 104000, 48534)

for ($i = 0; $i -lt 6; $i++) {
    $row = 16 + $i
    $ws.Cells.Item($row, 5).Value = $periods[$i]
    $ws.Cells.Item($row, 6).Value = $values[$i]
}

# 3. Remove the second worker (JAVIER ANTONIO HERRERA PALMERA / 1002299740) block entirely,
#    rows 22-27, shifting the remaining rows (signature block) up.
$ws.Rows("22:27").Delete()
